$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 447.08694
$ws.Range("I33").Value = 333.26315
$ws.Range("J33").Value = 987.75
$ws.Range("K33").Value = 333.26315
$ws.Range("L33").Value = 987.75
$ws.Range("M33").Value = -104.26315
$ws.Range("N33").Value = -1445.75
$ws.Range("H34").Value = 5120.1
$ws.Range("I34").Value = 5120.1
$ws.Range("K34").Value = 5120.1
$ws.Range("M34").Value = -4917.1
$ws.Range("H36").Value = 5120.1
$ws.Range("I36").Value = 5120.1
$ws.Range("K36").Value = 5120.1
$ws.Range("M36").Value = -4405.1
$ws.Range("H62").Value = 2499.7
$ws.Range("J62").Value = 3870
$ws.Range("L62").Value = 3870
$ws.Range("N62").Value = -5118
$ws.Range("H65").Value = 2499.7
$ws.Range("J65").Value = 3870
$ws.Range("L65").Value = 19350
$ws.Range("N65").Value = -25590
$ws.Range("H107").Value = 1803.4
$ws.Range("I107").Value = 1137.4166
$ws.Range("J107").Value = 4467.3335
$ws.Range("K107").Value = 1137.4166
$ws.Range("L107").Value = 4467.3335
$ws.Range("M107").Value = 782.5834
$ws.Range("N107").Value = -8307.333500000001
$ws.Range("H116").Value = 6591.143
$ws.Range("I116").Value = 6624.6665
$ws.Range("J116").Value = 6390
$ws.Range("K116").Value = 6624.6665
$ws.Range("L116").Value = 6390
$ws.Range("M116").Value = -3182.6665
$ws.Range("N116").Value = -13274
$ws.Range("H126").Value = 89131
$ws.Range("J126").Value = 89131
$ws.Range("L126").Value = 89131
$ws.Range("N126").Value = -99011
$ws.Range("H130").Value = 88749
$ws.Range("J130").Value = 88749
$ws.Range("L130").Value = 88749
$ws.Range("N130").Value = -98789
$ws.Range("H132").Value = 3044.2727
$ws.Range("I132").Value = 2855.9048
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 8567.714399999999
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -6037.714399999999
$ws.Range("N132").Value = -26060
$ws.Range("H141").Value = 5210.1816
$ws.Range("I141").Value = 3164.125
$ws.Range("K141").Value = 9492.375
$ws.Range("M141").Value = -4312.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6754.722
$ws.Range("I32").Value = 4105.6665
$ws.Range("K32").Value = 4105.6665
$ws.Range("M32").Value = -3818.6665
$ws.Range("H45").Value = 2065.5
$ws.Range("I45").Value = 1861.625
$ws.Range("K45").Value = 1861.625
$ws.Range("M45").Value = -1484.625
$ws.Range("H110").Value = 1510.8334
$ws.Range("I110").Value = 1344.5
$ws.Range("K110").Value = 1344.5
$ws.Range("M110").Value = 700.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H86").Value = 3821
$ws.Range("I86").Value = 4069.4
$ws.Range("K86").Value = 4069.4
$ws.Range("M86").Value = -2946.4
$ws.Range("H89").Value = 3821
$ws.Range("I89").Value = 4069.4
$ws.Range("K89").Value = 20347
$ws.Range("M89").Value = -14731
$ws.Range("H99").Value = 2740.75
$ws.Range("I99").Value = 2740.75
$ws.Range("K99").Value = 2740.75
$ws.Range("M99").Value = -1242.75
$ws.Range("H107").Value = 2304.6667
$ws.Range("I107").Value = 1511
$ws.Range("J107").Value = 3892
$ws.Range("K107").Value = 1511
$ws.Range("L107").Value = 3892
$ws.Range("M107").Value = 409
$ws.Range("N107").Value = -7732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.71429000000001
$ws.Range("I7").Value = 12.5
$ws.Range("J7").Value = 139
$ws.Range("K7").Value = 12.5
$ws.Range("L7").Value = 139
$ws.Range("M7").Value = 100.5
$ws.Range("N7").Value = -365
$ws.Range("H62").Value = 3048.125
$ws.Range("I62").Value = 3048.125
$ws.Range("K62").Value = 3048.125
$ws.Range("M62").Value = -2424.125
$ws.Range("H65").Value = 3048.125
$ws.Range("I65").Value = 3048.125
$ws.Range("K65").Value = 15240.625
$ws.Range("M65").Value = -12120.625
$ws.Range("H99").Value = 76999.28999999999
$ws.Range("I99").Value = 7332.6665
$ws.Range("J99").Value = 129249.25
$ws.Range("K99").Value = 7332.6665
$ws.Range("L99").Value = 129249.25
$ws.Range("M99").Value = -5834.6665
$ws.Range("N99").Value = -132245.25
$ws.Range("H126").Value = 76999.28999999999
$ws.Range("I126").Value = 7332.6665
$ws.Range("J126").Value = 129249.25
$ws.Range("K126").Value = 21997.9995
$ws.Range("L126").Value = 387747.75
$ws.Range("M126").Value = -19527.9995
$ws.Range("N126").Value = -392687.75
$ws.Range("H130").Value = 52877.5
$ws.Range("H134").Value = 3502.9333
$ws.Range("I134").Value = 2146
$ws.Range("K134").Value = 6438
$ws.Range("M134").Value = -3903

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 33682.668
$ws.Range("I11").Value = 50225
$ws.Range("K11").Value = 150675
$ws.Range("M11").Value = -150535
$ws.Range("H14").Value = 171.875
$ws.Range("I14").Value = 171.875
$ws.Range("K14").Value = 515.625
$ws.Range("M14").Value = -342.625
$ws.Range("H22").Value = 1270
$ws.Range("I22").Value = 1225
$ws.Range("K22").Value = 3675
$ws.Range("M22").Value = -3506
$ws.Range("H27").Value = 1270
$ws.Range("I27").Value = 1225
$ws.Range("K27").Value = 3675
$ws.Range("M27").Value = -3573
$ws.Range("H98").Value = 595.6667
$ws.Range("I98").Value = 697
$ws.Range("J98").Value = 393
$ws.Range("K98").Value = 2091
$ws.Range("L98").Value = 1179
$ws.Range("M98").Value = -593
$ws.Range("N98").Value = -4175
$ws.Range("H129").Value = 3990.6667
$ws.Range("I129").Value = 1330
$ws.Range("J129").Value = 4180.7144
$ws.Range("K129").Value = 3990
$ws.Range("L129").Value = 12542.1432
$ws.Range("N129").Value = -22542.1432
$ws.Range("M129").Value = 1010
$ws.Range("H131").Value = 23371.117
$ws.Range("I131").Value = 556254.5
$ws.Range("J131").Value = 1620.7755
$ws.Range("K131").Value = 1668763.5
$ws.Range("L131").Value = 4862.3265
$ws.Range("M131").Value = -1663723.5
$ws.Range("N131").Value = -14942.3265
$ws.Range("H136").Value = 7888.2
$ws.Range("I136").Value = 6702.5557
$ws.Range("K136").Value = 20107.6671
$ws.Range("M136").Value = -15007.6671
$ws.Range("H138").Value = 2456.3
$ws.Range("I138").Value = 1618.1111
$ws.Range("K138").Value = 4854.3333
$ws.Range("M138").Value = 285.6666999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8249.125
$ws.Range("I80").Value = 4666.4
$ws.Range("J80").Value = 14220.333
$ws.Range("K80").Value = 4666.4
$ws.Range("L80").Value = 14220.333
$ws.Range("M80").Value = -3668.4
$ws.Range("N80").Value = -16216.333
$ws.Range("H83").Value = 8249.125
$ws.Range("I83").Value = 4666.4
$ws.Range("J83").Value = 14220.333
$ws.Range("K83").Value = 23332
$ws.Range("L83").Value = 71101.66500000001
$ws.Range("M83").Value = -18340
$ws.Range("N83").Value = -81085.66500000001
$ws.Range("H102").Value = 4022
$ws.Range("I102").Value = 3629
$ws.Range("J102").Value = 4169.375
$ws.Range("K102").Value = 3629
$ws.Range("L102").Value = 4169.375
$ws.Range("M102").Value = -2007
$ws.Range("N102").Value = -7413.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1390.3334
$ws.Range("I22").Value = 653.8
$ws.Range("K22").Value = 653.8
$ws.Range("M22").Value = -358.8
$ws.Range("H27").Value = 1390.3334
$ws.Range("I27").Value = 653.8
$ws.Range("K27").Value = 653.8
$ws.Range("M27").Value = -546.8
$ws.Range("H35").Value = 1974.3334
$ws.Range("J35").Value = 1983
$ws.Range("L35").Value = 1983
$ws.Range("N35").Value = -2655
$ws.Range("H46").Value = 2566.9
$ws.Range("I46").Value = 1148.7778
$ws.Range("J46").Value = 3727.182
$ws.Range("K46").Value = 1148.7778
$ws.Range("L46").Value = 3727.182
$ws.Range("M46").Value = -960.7778000000001
$ws.Range("N46").Value = -4103.182
$ws.Range("H122").Value = 3843.3125
$ws.Range("I122").Value = 3847.5386
$ws.Range("J122").Value = 3825
$ws.Range("K122").Value = 11542.6158
$ws.Range("L122").Value = 11475
$ws.Range("M122").Value = -9092.6158
$ws.Range("N122").Value = -16375
$ws.Range("H127").Value = 51249.5
$ws.Range("J127").Value = 51249.5
$ws.Range("L127").Value = 51249.5
$ws.Range("N127").Value = -61169.5
$ws.Range("H130").Value = 56985.6
$ws.Range("J130").Value = 56985.6
$ws.Range("L130").Value = 56985.6
$ws.Range("N130").Value = -67025.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 42808.855
$ws.Range("J112").Value = 42808.855
$ws.Range("L112").Value = 42808.855
$ws.Range("N112").Value = -45762.855
$ws.Range("H126").Value = 3934.2144
$ws.Range("I126").Value = 4054.2
$ws.Range("J126").Value = 3634.25
$ws.Range("K126").Value = 12162.6
$ws.Range("L126").Value = 10902.75
$ws.Range("M126").Value = -9692.599999999999
$ws.Range("N126").Value = -15842.75
